# Packaging_Data worksheet update (release v1.0, Sept 19)
#
# 1. Drop the now-unused helper column L (was entirely empty, styled only).
# 2. Update a handful of dimension inputs that moved from computed
#    "n * 25.4" conversions to directly-entered millimetre values.
# 3. Tidy the header row's height now that column L (and its wrap-text
#    formatting) is gone.
# 4. Restore the active cell/selection to where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Remove column L entirely -------------------------------------------
$ws.Columns("L").Delete()

# Re-fit the header row height now that the L column (and its formatting)
# is gone - drops the stale explicit ht="45".
$ws.Rows("1").AutoFit()

# --- 2. Data edits -----------------------------------------------------------

# Pack_Angle: Height (mm) 75 -> 5
$ws.Range("D2").Value = 5

# Pack_1x4x19_Dun: Width/Height now entered directly instead of via "n*25.4"
$ws.Range("C5").Value = 89
$ws.Range("D5").Value = 19

# Pack_2x3x19_Dun: Width/Height now entered directly instead of via "n*25.4"
$ws.Range("C6").Value = 64
$ws.Range("D6").Value = 38

# COMMON LUMBER (Pack row 7): Width/Height now entered directly
$ws.Range("C7").Value = 89
$ws.Range("D7").Value = 19

# --- 3. Restore cursor/selection position -----------------------------------
$ws.Range("D20").Select() | Out-Null
